# Applies the "Added a new Time Difference cost function. Not as powerful"
# commit to the TDoA sensitivity-analysis sheet.
#
# Summary of the change:
#  - Rows 18-24, 26-27: recomputed cost-function figures (columns B-H)
#  - Row 20 & 24: the "Largest Angle" column (G) becomes the literal text
#    "NaN" (same text already used in G22) and H becomes 0
#  - Row 25: recomputed figures AND the whole row gets a yellow highlight
#    (new cost function row)
#  - Row 28 (new row): another yellow-highlighted row with the new
#    cost-function figures, column A left blank
#  - Selection moves to C29 or roughly that area

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Updated numeric data for existing rows 18-27 (columns B..H)
#    NaN is written as literal text "NaN" so Excel stores/text it
#    exactly like the pre-existing G22 cell.
# ---------------------------------------------------------------------
$rows = @(
  @{R=18; B=0.58796296296296302; C=6.5041224460093501;  D=28.214845749269699;  E=3.9240778458147201; F=2.1766512126122302; G=25;    H=60}
  @{R=19; B=0.62654320987654299; C=5.3975741344306396;  D=35.741776044801703;  E=3.1007782769840899; F=1.2860353659228601; G=25;    H=65}
  @{R=20; B=0;                   C=183.186111054778;    D=92.223100532066297;  E=14.6436579027897;   F=6.2363243441009297; G="NaN"; H=0}
  @{R=21; B=0.50462962962962998; C=9.6073890320557496;  D=71.090035901500997;  E=4.3709213736600798; F=1.82803638148556;   G=30;    H=55}
  @{R=22; B=0.296296296296296;   C=30.140142176298699;  D=113.34225657166201;  E=6.64992657391634;   F=3.0793490493221198; G=$null; H=$null}
  @{R=23; B=0.67592592592592604; C=4.2749364100667302;  D=18.156706614373402;  E=2.8269370789892898; F=1.31865104338513;   G=20;    H=65}
  @{R=24; B=0.31018518518518501; C=30.785045539260899;  D=121.68714370284199;  E=6.0031108031544198; F=3.1910591821854202; G="NaN"; H=0}
  @{R=25; B=0.68364197530864201; C=3.78048484491296;    D=15.7749639368175;    E=2.7121008342665802; F=1.3893084645587701; G=20;    H=70}
  @{R=26; B=0.66358024691357997; C=4.6177686190237104;  D=14.763947564591399;  E=2.9205206881378398; F=1.3392072681690901; G=20;    H=$null}
  @{R=27; B=0.60030864197530898; C=6.5492757791657796;  D=44.445453165045699;  E=3.89762276430173;   F=1.7636480260990399; G=20;    H=55}
)

foreach ($d in $rows) {
  $r = $d.R
  $ws.Range("B$r").Value = $d.B
  $ws.Range("C$r").Value = $d.C
  $ws.Range("D$r").Value = $d.D
  $ws.Range("E$r").Value = $d.E
  $ws.Range("F$r").Value = $d.F
  if ($null -ne $d.G) { $ws.Range("G$r").Value = $d.G }
  if ($null -ne $d.H) { $ws.Range("H$r").Value = $d.H }
}

# ---------------------------------------------------------------------
# 2. Row 25 becomes the new highlighted "cost function" row - apply the
#    yellow fill to A25:H25 (keeps each column's own number format).
# ---------------------------------------------------------------------
$ws.Range("A25:H25").Interior.Color = 65535

# ---------------------------------------------------------------------
# 3. Brand new row 28 - another highlighted cost-function row, inserted
#    right above the pre-existing blank spacer row 29. Column A is left
#    empty (only the yellow fill carries over).
# ---------------------------------------------------------------------
$ws.Range("B28").Value = 0.657407407407407
$ws.Range("C28").Value = 4.0083047628767803
$ws.Range("D28").Value = 23.027850607104199
$ws.Range("E28").Value = 2.7258275260840401
$ws.Range("F28").Value = 1.5271403356623601
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = 70

# Give row 28 the same number formats as row 25 so the new cellXfs match
# (percentage, 0.0, integer, general) before/while applying the fill.
$ws.Range("A28").NumberFormat = "General"
$ws.Range("B28").NumberFormat = $ws.Range("B25").NumberFormat
$ws.Range("C28:F28").NumberFormat = $ws.Range("C25").NumberFormat
$ws.Range("G28").NumberFormat = $ws.Range("G25").NumberFormat
$ws.Range("H28").NumberFormat = "General"

$ws.Range("A28:H28").Interior.Color = 65535

# ---------------------------------------------------------------------
# 4. Move the active selection roughly where the author left it.
# ---------------------------------------------------------------------
$ws.Range("C29").Select() | Out-Null
